# Apply the update described by the diff:
#  - Append rows 122-130 to the "Orders" sheet (sheet1)
#  - Update the concatenated "TotalNumber" text in cell G2 of the "Summary" sheet (sheet2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Data for the new rows: row number, PackageID (A), FlowerName (C), Number (F)
$newRows = @(
    @{ Row = 122; A = $null; C = "80_冰清玉洁_undefined_Gerbera L._20stems"; F = "10" },
    @{ Row = 123; A = $null; C = "77_珍爱mini_undefined_Gerbera L._20stems"; F = "10" },
    @{ Row = 124; A = $null; C = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"; F = "8" },
    @{ Row = 125; A = $null; C = "631_吸色康乃馨宝蓝_tinted blue_undefined_20stems"; F = "10" },
    @{ Row = 126; A = $null; C = "631_吸色康乃馨宝蓝_tinted blue_undefined_20stems"; F = "10" },
    @{ Row = 127; A = "4";   C = "669_大丽花 红_undefined_undefined_5stems"; F = "10" },
    @{ Row = 128; A = $null; C = "420_松虫草QQ糖_scabiosa white pink_undefined_1bunch"; F = "15" },
    @{ Row = 129; A = $null; C = "418_松虫草白_scabiosa white_undefined_1bunch"; F = "10" },
    @{ Row = 130; A = $null; C = "648_洋牡丹河内_undefined_undefined_1bunch"; F = "5" }
)

foreach ($r in $newRows) {
    if ($r.A -ne $null) {
        $cellA = $ws.Cells.Item($r.Row, 1)
        $cellA.NumberFormat = "@"
        $cellA.Value = $r.A
    }

    $cellC = $ws.Cells.Item($r.Row, 3)
    $cellC.Value = $r.C

    # Force column F to be stored as text (matches the source workbook, where
    # every "Number" value is kept as text, not a real number).
    $cellF = $ws.Cells.Item($r.Row, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $r.F
}

# Update the Summary sheet's concatenated TotalNumber string (G2) to include
# the newly-added rows' Number values appended to the existing text.
$summary = $wb.Worksheets.Item("Summary")
$g2 = $summary.Cells.Item(2, 7)
$g2.NumberFormat = "@"
$g2.Value = "055155552510652566555525321515822555510555551255156558101576510612610551051510555510158105151051541156111387865775125551281055451565551010352.515202055105121141410152020201010810101015105"
